$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: add note in column B
$ws.Range("B18").Value = "Setting up YOLO"

# Row 19: new date entry (formatted like the rows above it)
$ws.Range("A19").Value = 43774
$ws.Range("A19").NumberFormat = "MM/DD/YY"

# Row 20: new date entry (formatted like the rows above it)
$ws.Range("A20").Value = 43775
$ws.Range("A20").NumberFormat = "MM/DD/YY"

# Update selection to match the new active cell
$ws.Range("B19").Select()
